# Added Panel Accessories Test Data For Spain/Turkey/Hungary market
#
# For each of the three market sheets (Turkey, Hungary, Spain) two new
# accessory rows ("MX-BBX" / "MX-DPBX") are inserted into the accessories
# block at the bottom of the sheet, just above the row(s) that already
# follow the last "real" accessory row. The new rows re-use the
# formatting (border style) of the row directly above the insertion
# point, and the new cell values come from the shared strings already
# present in the workbook (MX-BBX / MX-DPBX).

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

function Add-AccessoryRows($SheetName, $InsertBeforeRow) {
    $ws = $wb.Worksheets.Item($SheetName)

    # Row immediately above the insertion point donates its formatting
    # (style index) to the two freshly inserted rows.
    $formatSourceRow = $InsertBeforeRow - 1
    $secondRow = $InsertBeforeRow + 1

    $ws.Rows("$InsertBeforeRow`:$secondRow").Insert()

    $ws.Range("A$formatSourceRow").Copy()
    $ws.Range("A$InsertBeforeRow`:A$secondRow").PasteSpecial($xlPasteFormats)
    $excel.CutCopyMode = $false

    $ws.Range("A$InsertBeforeRow").Value = "MX-BBX"
    $ws.Range("A$secondRow").Value = "MX-DPBX"
}

# Turkey: new rows land right above the existing "PR1D2-Unmonitored" row (13).
Add-AccessoryRows "Turkey" 13
$wb.Worksheets.Item("Turkey").Range("A13:A14").Select()

# Hungary: sheet has no "PR1D2-Unmonitored" row, so the new rows land just
# above the trailing "Wg"/"Accessories" rows (row 14).
Add-AccessoryRows "Hungary" 14
$wb.Worksheets.Item("Hungary").Range("A14:A15").Select()

# Spain: same layout as Turkey - new rows land above "PR1D2-Unmonitored" (13).
Add-AccessoryRows "Spain" 13

# Poland was the previously-active sheet; its selection moved to A9:A10
# before the user switched over to Spain.
$wb.Worksheets.Item("Poland").Range("A9:A10").Select()

# Spain ends up as the active sheet with B12 selected.
$wsSpain = $wb.Worksheets.Item("Spain")
$wsSpain.Activate()
$wsSpain.Range("B12").Select()
